# Edit script: implements the diff changes against before.docx
$d = $word.ActiveDocument

function Replace-ParaRuns($Index, $PPrXml, $RunsXml) {
    $p = $d.Paragraphs($Index)
    $start = $p.Range.Start
    $end = $p.Range.End - 1
    $r = $d.Range($start, $end)
    $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + $PPrXml + $RunsXml + "</w:p>"
    $r.InsertXML($xml)
}

$pprNormal0 = "<w:pPr><w:pStyle w:val='Normal'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr><w:rPr/></w:pPr>"
$pprNormal1 = "<w:pPr><w:pStyle w:val='Normal'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr><w:rPr/></w:pPr>"

# Paragraph 1: split into two runs, removing "except for the final API endpoint," tail
$runs1 = "<w:r><w:rPr/><w:t xml:space='preserve'>The project is passing all functional test cases. The project successfully implements all API endpoints. </w:t></w:r>" + "<w:r><w:rPr/><w:t>But the Search class is not unit tested</w:t></w:r>"
Replace-ParaRuns 1 $pprNormal0 $runs1

# Paragraph 5: "Lines of code in unit tests: 539" -> two runs "Lines of code in unit tests: " + "729"
$runs5 = "<w:r><w:rPr/><w:t xml:space='preserve'>Lines of code in unit tests: </w:t></w:r>" + "<w:r><w:rPr/><w:t>729</w:t></w:r>"
Replace-ParaRuns 5 $pprNormal0 $runs5

# Paragraph 6: "Unit test coverage measured by tap-istanbul: 98.62%" -> five runs
$runs6 = "<w:r><w:rPr/><w:t>Unit test coverage measured by tap-istanbul: 9</w:t></w:r>" + "<w:r><w:rPr/><w:t>2</w:t></w:r>" + "<w:r><w:rPr/><w:t>.</w:t></w:r>" + "<w:r><w:rPr/><w:t>49</w:t></w:r>" + "<w:r><w:rPr/><w:t>%</w:t></w:r>"
Replace-ParaRuns 6 $pprNormal0 $runs6

# Paragraph 8: "7 hours to completion" -> "32" + " hours to completion"
$runs8 = "<w:r><w:rPr/><w:t>32</w:t></w:r>" + "<w:r><w:rPr/><w:t xml:space='preserve'> hours to completion</w:t></w:r>"
Replace-ParaRuns 8 $pprNormal0 $runs8

# Paragraph 9: "1 hour was spent preparing the submission" -> "2" + " hour was spent preparing the submission"
$runs9 = "<w:r><w:rPr/><w:t>2</w:t></w:r>" + "<w:r><w:rPr/><w:t xml:space='preserve'> hour was spent preparing the submission</w:t></w:r>"
Replace-ParaRuns 9 $pprNormal0 $runs9

# Paragraph 11: "Achieving above 80% ... till 98% code coverage ..." -> three runs
$runs11 = "<w:r><w:rPr/><w:t>Achieving above 80% unit test coverage: I was able to get it above and all the way till 9</w:t></w:r>" + "<w:r><w:rPr/><w:t>2</w:t></w:r>" + "<w:r><w:rPr/><w:t>% code coverage by testing all my base classes thoroughly</w:t></w:r>"
Replace-ParaRuns 11 $pprNormal1 $runs11

Write-Host "Text edits complete"
for ($i=1; $i -le $d.Paragraphs.Count; $i++) {
  Write-Host "$i : [$($d.Paragraphs($i).Range.Text)]"
}
